# Generate Report for Handoff
# Renames the handed-off file ids, refreshes status/dates, and clears the
# now-redundant "Latest Target File" columns (I/J) on the per-locale sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "bf77464a-9512-4cdc-af40-f1a4f44a6b8a.md"
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"
$ws1.Range("G2").Value = "2016-08-19 21:05:49"

$ws1.Range("A3").Value = "ffff81406dd3-3959-4370-857a-33bb9c4072f5.md"
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-19 21:05:49"

# Hyperlink display text changes (B2/B3) while keeping the original targets.
$hls1 = $ws1.Hyperlinks
$addrs1 = @()
foreach ($hl in $hls1) {
    $addrs1 += $hl.Address
}
$hls1.Delete()
$hls1.Add($ws1.Range("B2"), $addrs1[0], "", "", "e2e\bf77464a-9512-4cdc-af40-f1a4f44a6b8a.md")
$hls1.Add($ws1.Range("B3"), $addrs1[1], "", "", "e2e\ffff81406dd3-3959-4370-857a-33bb9c4072f5.md")


# Column width tweaks on the Overview sheet (E & F narrower).
$ws1.Columns.Item(5).ColumnWidth = 16.333333333333332
$ws1.Columns.Item(6).ColumnWidth = 16.333333333333332

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "bf77464a-9512-4cdc-af40-f1a4f44a6b8a.md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("G2").Value = "bf77464a-9512-4cdc-af40-f1a4f44a6b8a.68e97bedd721f2a9b3f82b9feaf711f519f63669.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-19 21:05:45"
$ws2.Range("I2").Value = ""
$ws2.Range("J2").Value = ""
$ws2.Range("K2").Value = "0001-01-01 00:00:00"

$ws2.Range("A3").Value = "ffff81406dd3-3959-4370-857a-33bb9c4072f5.md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "bf77464a-9512-4cdc-af40-f1a4f44a6b8a.68e97bedd721f2a9b3f82b9feaf711f519f63669.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-19 21:05:45"
$ws2.Range("I3").Value = ""
$ws2.Range("J3").Value = ""
$ws2.Range("K3").Value = "0001-01-01 00:00:00"

# I2/I3 no longer carry the hyperlink style (now plain cells).
$ws2.Range("I2").Style = "Normal"
$ws2.Range("I3").Style = "Normal"

# Keep only the A2/A3 hyperlinks (drop I2/I3), preserving original targets.
$hls2 = $ws2.Hyperlinks
$addrs2 = @()
$refs2 = @()
foreach ($hl in $hls2) {
    $addrs2 += $hl.Address
    $refs2 += $hl.Range.Address()
}
$hls2.Delete()
for ($i = 0; $i -lt $refs2.Count; $i++) {
    if ($refs2[$i] -eq "`$A`$2") {
        $hls2.Add($ws2.Range("A2"), $addrs2[$i], "", "", "bf77464a-9512-4cdc-af40-f1a4f44a6b8a.md")
    }
    elseif ($refs2[$i] -eq "`$A`$3") {
        $hls2.Add($ws2.Range("A3"), $addrs2[$i], "", "", "ffff81406dd3-3959-4370-857a-33bb9c4072f5.md")
    }
}

# Column width tweaks on zh-cn sheet (C, I, J).
$ws2.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws2.Columns.Item(9).ColumnWidth = 17.833333333333332
$ws2.Columns.Item(10).ColumnWidth = 20.833333333333332

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "bf77464a-9512-4cdc-af40-f1a4f44a6b8a.md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("G2").Value = "bf77464a-9512-4cdc-af40-f1a4f44a6b8a.68e97bedd721f2a9b3f82b9feaf711f519f63669.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-19 21:05:49"
$ws3.Range("I2").Value = ""
$ws3.Range("J2").Value = ""
$ws3.Range("K2").Value = "0001-01-01 00:00:00"

$ws3.Range("A3").Value = "ffff81406dd3-3959-4370-857a-33bb9c4072f5.md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "bf77464a-9512-4cdc-af40-f1a4f44a6b8a.68e97bedd721f2a9b3f82b9feaf711f519f63669.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-19 21:05:49"
$ws3.Range("I3").Value = ""
$ws3.Range("J3").Value = ""
$ws3.Range("K3").Value = "0001-01-01 00:00:00"

# I2/I3 no longer carry the hyperlink style (now plain cells).
$ws3.Range("I2").Style = "Normal"
$ws3.Range("I3").Style = "Normal"

# Keep only the A2/A3 hyperlinks (drop I2/I3), preserving original targets.
$hls3 = $ws3.Hyperlinks
$addrs3 = @()
$refs3 = @()
foreach ($hl in $hls3) {
    $addrs3 += $hl.Address
    $refs3 += $hl.Range.Address()
}
$hls3.Delete()
for ($i = 0; $i -lt $refs3.Count; $i++) {
    if ($refs3[$i] -eq "`$A`$2") {
        $hls3.Add($ws3.Range("A2"), $addrs3[$i], "", "", "bf77464a-9512-4cdc-af40-f1a4f44a6b8a.md")
    }
    elseif ($refs3[$i] -eq "`$A`$3") {
        $hls3.Add($ws3.Range("A3"), $addrs3[$i], "", "", "ffff81406dd3-3959-4370-857a-33bb9c4072f5.md")
    }
}

# Column width tweaks on de-de sheet (C, I, J).
$ws3.Columns.Item(3).ColumnWidth = 16.333333333333332
$ws3.Columns.Item(9).ColumnWidth = 17.833333333333332
$ws3.Columns.Item(10).ColumnWidth = 20.833333333333332
